# Daily refresh of the cryptos price list (GitHub Actions bot).
# Updates Price (D) / Volume(1h) (E) for most rows, and for rows 49-50
# the Stellar / WhiteBITCoin entries swap ranking order (with new values).
# A leading "'" forces a handful of Price cells that would otherwise be
# read back by Excel as a trimmed number (e.g. "1.00" -> 1) to stay text,
# matching the original exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.728.22"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "2.906.71"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "528.81"
$ws.Range("E5").Value = "  -1.50%  "

$ws.Range("D6").Value = "144.84"
$ws.Range("E6").Value = "  -4.29%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  -1.00%  "

$ws.Range("D9").Value = "2.912.87"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  -2.66%  "

$ws.Range("D11").Value = "'6.00"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").Value = "0.364"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").Value = "3.416.38"
$ws.Range("E13").Value = "  -1.53%  "

$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").Value = "60.701.35"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("D17").Value = "2.916.97"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("E18").Value = "  -2.18%  "

$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").Value = "11.67"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").Value = "361.73"
$ws.Range("E21").Value = "  -5.02%  "

$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").Value = "5.69"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").Value = "64.73"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "0.455"
$ws.Range("E26").Value = "  -2.18%  "

$ws.Range("D27").Value = "0.182"
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").Value = "7.83"
$ws.Range("E29").Value = "  -4.80%  "

$ws.Range("D30").Value = "0.0₃0866"
$ws.Range("E30").Value = "  -5.60%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").Value = "19.73"
$ws.Range("E33").Value = "  -2.54%  "

$ws.Range("D34").Value = "150.46"
$ws.Range("E34").Value = "  -4.84%  "

$ws.Range("D35").Value = "'4.40"
$ws.Range("E35").Value = "  -3.26%  "

$ws.Range("D36").Value = "5.58"
$ws.Range("E36").Value = "  -6.13%  "

$ws.Range("E37").Value = "  -4.29%  "

$ws.Range("E38").Value = "  -4.49%  "

$ws.Range("D39").Value = "37.73"
$ws.Range("E39").Value = "  +2.65%  "

$ws.Range("E40").Value = "  -2.13%  "

$ws.Range("E41").Value = "  -3.73%  "

$ws.Range("D42").Value = "2.287.18"
$ws.Range("E42").Value = "  -4.78%  "

$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").Value = "0.0583"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").Value = "20.76"
$ws.Range("E45").Value = "  -5.71%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("E47").Value = "  +4.01%  "

$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "10.34"
$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.0926"
$ws.Range("E50").Value = "  -1.67%  "

$ws.Range("D51").Value = "251.59"
$ws.Range("E51").Value = "  -4.44%  "
